$d = $word.ActiveDocument

$d.Content.Find.Execute("31×43=", $true, $false, $false, $false, $false, $true, 1, $false, "59×48=", 2) | Out-Null
$d.Content.Find.Execute("30×42=", $true, $false, $false, $false, $false, $true, 1, $false, "15×28=", 2) | Out-Null
$d.Content.Find.Execute("40×78=", $true, $false, $false, $false, $false, $true, 1, $false, "86×24=", 2) | Out-Null
$d.Content.Find.Execute("58×40=", $true, $false, $false, $false, $false, $true, 1, $false, "52×26=", 2) | Out-Null
$d.Content.Find.Execute("83×25=", $true, $false, $false, $false, $false, $true, 1, $false, "13×37=", 2) | Out-Null
$d.Content.Find.Execute("29×27=", $true, $false, $false, $false, $false, $true, 1, $false, "89×65=", 2) | Out-Null
$d.Content.Find.Execute("90×42=", $true, $false, $false, $false, $false, $true, 1, $false, "40×34=", 2) | Out-Null
$d.Content.Find.Execute("14×33=", $true, $false, $false, $false, $false, $true, 1, $false, "29×66=", 2) | Out-Null
$d.Content.Find.Execute("48×86=", $true, $false, $false, $false, $false, $true, 1, $false, "41×14=", 2) | Out-Null
$d.Content.Find.Execute("50×98=", $true, $false, $false, $false, $false, $true, 1, $false, "64×38=", 2) | Out-Null
$d.Content.Find.Execute("51×15=", $true, $false, $false, $false, $false, $true, 1, $false, "70×21=", 2) | Out-Null
$d.Content.Find.Execute("36×79=", $true, $false, $false, $false, $false, $true, 1, $false, "46×54=", 2) | Out-Null
$d.Content.Find.Execute("70×72=", $true, $false, $false, $false, $false, $true, 1, $false, "47×88=", 2) | Out-Null
$d.Content.Find.Execute("74×11=", $true, $false, $false, $false, $false, $true, 1, $false, "67×19=", 2) | Out-Null
$d.Content.Find.Execute("44×25=", $true, $false, $false, $false, $false, $true, 1, $false, "24×28=", 2) | Out-Null
$d.Content.Find.Execute("52×74=", $true, $false, $false, $false, $false, $true, 1, $false, "74×39=", 2) | Out-Null
$d.Content.Find.Execute("38×77=", $true, $false, $false, $false, $false, $true, 1, $false, "94×39=", 2) | Out-Null
$d.Content.Find.Execute("79×90=", $true, $false, $false, $false, $false, $true, 1, $false, "18×33=", 2) | Out-Null
$d.Content.Find.Execute("11×40=", $true, $false, $false, $false, $false, $true, 1, $false, "60×79=", 2) | Out-Null
$d.Content.Find.Execute("38×86=", $true, $false, $false, $false, $false, $true, 1, $false, "59×15=", 2) | Out-Null
$d.Content.Find.Execute("19×53=", $true, $false, $false, $false, $false, $true, 1, $false, "75×60=", 2) | Out-Null
$d.Content.Find.Execute("13×28=", $true, $false, $false, $false, $false, $true, 1, $false, "32×51=", 2) | Out-Null
$d.Content.Find.Execute("27×94=", $true, $false, $false, $false, $false, $true, 1, $false, "49×41=", 2) | Out-Null
$d.Content.Find.Execute("74×84=", $true, $false, $false, $false, $false, $true, 1, $false, "87×37=", 2) | Out-Null
$d.Content.Find.Execute("37×71=", $true, $false, $false, $false, $false, $true, 1, $false, "69×47=", 2) | Out-Null
